$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = 0.0001625819899935208
$ws.Range("D2").Value = 0.5885782825730005
$ws.Range("E2").Value = 0.6492109438477345
$ws.Range("B3").Value = [double]"2.183980475909258E-10"
$ws.Range("C3").Value = 0.007919668242014082
$ws.Range("D3").Value = 0.5501068836496137
$ws.Range("E3").Value = 0.5498836806399565
$ws.Range("B4").Value = [double]"3.409259119931335E-12"
$ws.Range("C4").Value = 0.00715402768698076
$ws.Range("D4").Value = 0.3927863988361095
$ws.Range("E4").Value = 0.460358406415663
$ws.Range("C5").Value = [double]"1.587093575472108E-07"
$ws.Range("D5").Value = 0.02009765232936478
$ws.Range("E5").Value = 0.03658129441505157

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = 0.0001842239112421771
$ws.Range("D2").Value = 0.6669262277582069
$ws.Range("E2").Value = 0.7356299384797864
$ws.Range("B3").Value = [double]"2.474698614308933E-10"
$ws.Range("C3").Value = 0.008973886094901615
$ws.Range("D3").Value = 0.6233337512427775
$ws.Range("E3").Value = 0.6230808368120845
$ws.Range("B4").Value = [double]"3.863078865850037E-12"
$ws.Range("C4").Value = 0.008106328146696585
$ws.Range("D4").Value = 0.4450717246061602
$ws.Range("E4").Value = 0.5216385050182306
$ws.Range("C5").Value = [double]"1.798357776236215E-07"
$ws.Range("D5").Value = 0.02277292902521752
$ws.Range("E5").Value = 0.0414507728421276

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = [double]"6.556994183442194E-06"
$ws.Range("C2").Value = 0.005083743701291375
$ws.Range("D2").Value = 1.07580004858252
$ws.Range("E2").Value = 0.9697180417558001
$ws.Range("B3").Value = [double]"4.458740202226576E-05"
$ws.Range("C3").Value = 0.01839468255529511
$ws.Range("D3").Value = 0.7727021730903773
$ws.Range("E3").Value = 0.6805195557530022
$ws.Range("B4").Value = 0.0001322419804487943
$ws.Range("C4").Value = 0.004894520845755111
$ws.Range("D4").Value = 0.5532047049309021
$ws.Range("E4").Value = 0.6055420661490584
$ws.Range("B5").Value = [double]"4.153988879660067E-05"
$ws.Range("C5").Value = 0.01081418683615371
$ws.Range("D5").Value = 1.024873148138431
$ws.Range("E5").Value = 0.7122850722825392

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = [double]"4.966528370999151E-05"
$ws.Range("C2").Value = 0.003462735458016471
$ws.Range("D2").Value = 1.385766416475183
$ws.Range("E2").Value = 1.39053799804532
$ws.Range("B3").Value = [double]"5.30035999530297E-05"
$ws.Range("C3").Value = 0.01164057132748931
$ws.Range("D3").Value = 0.6708009055267363
$ws.Range("E3").Value = 0.648446477114007
$ws.Range("B4").Value = 0.0003397047964529607
$ws.Range("C4").Value = 0.003246107553390979
$ws.Range("D4").Value = 0.6270289992800223
$ws.Range("E4").Value = 0.7025531601665833
$ws.Range("B5").Value = 0.0001823860200208515
$ws.Range("C5").Value = 0.004124207012744106
$ws.Range("D5").Value = 1.240009647204811
$ws.Range("E5").Value = 0.9771097180398071
